# Auto-generated Excel COM-interop script applying the Sheets profit recalculation
# (currentAveragePrice / LevePrice / LeveProfit columns) described in the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 15555.556
$ws.Range("J17").Value = 15555.556
$ws.Range("L17").Value = 46666.66800000001
$ws.Range("N17").Value = -47002.66800000001

# Row 38
$ws.Range("H38").Value = 1335.8667
$ws.Range("I38").Value = 1335.8667
$ws.Range("K38").Value = 4007.6001
$ws.Range("M38").Value = -3635.6001

# Row 70
$ws.Range("H70").Value = 1972.75
$ws.Range("I70").Value = 1899
$ws.Range("J70").Value = 1997.3334
$ws.Range("K70").Value = 5697
$ws.Range("L70").Value = 5992.0002
$ws.Range("M70").Value = -5427
$ws.Range("N70").Value = -6532.0002

# Row 73
$ws.Range("H73").Value = 1972.75
$ws.Range("I73").Value = 1899
$ws.Range("J73").Value = 1997.3334
$ws.Range("K73").Value = 5697
$ws.Range("L73").Value = 5992.0002
$ws.Range("M73").Value = -4761
$ws.Range("N73").Value = -7864.0002

# Row 116
$ws.Range("H116").Value = 8319.385
$ws.Range("I116").Value = 7914.75
$ws.Range("K116").Value = 7914.75
$ws.Range("M116").Value = -4472.75

# Row 137
$ws.Range("H137").Value = 2945.6667
$ws.Range("I137").Value = 2311
$ws.Range("J137").Value = 4215
$ws.Range("K137").Value = 6933
$ws.Range("L137").Value = 12645
$ws.Range("M137").Value = -4383
$ws.Range("N137").Value = -17745

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1202.5
$ws.Range("I45").Value = 1202.5
$ws.Range("K45").Value = 1202.5
$ws.Range("M45").Value = -825.5

# Row 61
$ws.Range("H61").Value = 5084.4116
$ws.Range("J61").Value = 7998.778
$ws.Range("L61").Value = 7998.778
$ws.Range("N61").Value = -8422.778

# Row 74
$ws.Range("H74").Value = 2848.842
$ws.Range("J74").Value = 4110.143
$ws.Range("L74").Value = 4110.143
$ws.Range("N74").Value = -5858.143

# Row 77
$ws.Range("H77").Value = 2848.842
$ws.Range("J77").Value = 4110.143
$ws.Range("L77").Value = 20550.715
$ws.Range("N77").Value = -29286.715

# Row 122
$ws.Range("H122").Value = 2003.9
$ws.Range("I122").Value = 1581.6666
$ws.Range("J122").Value = 2637.25
$ws.Range("K122").Value = 4744.9998
$ws.Range("L122").Value = 7911.75
$ws.Range("M122").Value = -2294.9998
$ws.Range("N122").Value = -12811.75

# Row 130
$ws.Range("H130").Value = 50782.25
$ws.Range("J130").Value = 50782.25
$ws.Range("L130").Value = 50782.25
$ws.Range("N130").Value = -60822.25

# Row 136
$ws.Range("H136").Value = 5084.4116
$ws.Range("J136").Value = 7998.778
$ws.Range("L136").Value = 23996.334
$ws.Range("N136").Value = -29096.334

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 977.3333
$ws.Range("I20").Value = 1047
$ws.Range("K20").Value = 1047
$ws.Range("M20").Value = -800

# Row 134
$ws.Range("H134").Value = 1870.6154
$ws.Range("I134").Value = 1870.6154
$ws.Range("K134").Value = 5611.8462
$ws.Range("M134").Value = -3076.8462

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1713
$ws.Range("I16").Value = 1209.5385
$ws.Range("K16").Value = 1209.5385
$ws.Range("M16").Value = -922.5385000000001

# Row 22
$ws.Range("H22").Value = 43247.332
$ws.Range("I22").Value = 2746.3333
$ws.Range("J22").Value = 83748.336
$ws.Range("K22").Value = 2746.3333
$ws.Range("L22").Value = 83748.336
$ws.Range("M22").Value = -2396.3333
$ws.Range("N22").Value = -84448.336

# Row 31
$ws.Range("H31").Value = 3494.4443
$ws.Range("I31").Value = 2673.75
$ws.Range("K31").Value = 2673.75
$ws.Range("M31").Value = -2378.75

# Row 34
$ws.Range("H34").Value = 3494.4443
$ws.Range("I34").Value = 2673.75
$ws.Range("K34").Value = 2673.75
$ws.Range("M34").Value = -2471.75

# Row 58
$ws.Range("H58").Value = 1946.4615
$ws.Range("I58").Value = 1892.7142
$ws.Range("K58").Value = 1892.7142
$ws.Range("M58").Value = -1689.7142

# Row 107
$ws.Range("H107").Value = 909.9
$ws.Range("I107").Value = 531
$ws.Range("J107").Value = 1478.25
$ws.Range("K107").Value = 531
$ws.Range("L107").Value = 1478.25
$ws.Range("M107").Value = 1389
$ws.Range("N107").Value = -5318.25

# Row 113
$ws.Range("H113").Value = 1713
$ws.Range("I113").Value = 1209.5385
$ws.Range("K113").Value = 1209.5385
$ws.Range("M113").Value = 960.4614999999999

# Row 122
$ws.Range("H122").Value = 2027.1
$ws.Range("I122").Value = 1704.75
$ws.Range("J122").Value = 2242
$ws.Range("K122").Value = 5114.25
$ws.Range("L122").Value = 6726
$ws.Range("M122").Value = -2664.25
$ws.Range("N122").Value = -11626

# Row 136
$ws.Range("H136").Value = 1946.4615
$ws.Range("I136").Value = 1892.7142
$ws.Range("K136").Value = 5678.142599999999
$ws.Range("M136").Value = -3128.142599999999

$ws = $wb.Worksheets.Item("CUL")
# Row 94
$ws.Range("H94").Value = 10716
$ws.Range("I94").Value = 2972.5715
$ws.Range("K94").Value = 8917.7145
$ws.Range("M94").Value = -8241.7145

# Row 131
$ws.Range("H131").Value = 2823.2
$ws.Range("I131").Value = 1933
$ws.Range("J131").Value = 3045.75
$ws.Range("K131").Value = 5799
$ws.Range("L131").Value = 9137.25
$ws.Range("M131").Value = -759
$ws.Range("N131").Value = -19217.25

# Row 139
$ws.Range("H139").Value = 3368.6
$ws.Range("I139").Value = 3368.6
$ws.Range("K139").Value = 10105.8
$ws.Range("M139").Value = -4965.799999999999

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 2764.7778
$ws.Range("I80").Value = 1462.6
$ws.Range("J80").Value = 3265.6155
$ws.Range("K80").Value = 1462.6
$ws.Range("L80").Value = 3265.6155
$ws.Range("M80").Value = -464.5999999999999
$ws.Range("N80").Value = -5261.6155

# Row 83
$ws.Range("H83").Value = 2764.7778
$ws.Range("I83").Value = 1462.6
$ws.Range("J83").Value = 3265.6155
$ws.Range("K83").Value = 7313
$ws.Range("L83").Value = 16328.0775
$ws.Range("M83").Value = -2321
$ws.Range("N83").Value = -26312.0775

# Row 113
$ws.Range("H113").Value = 1356.5714
$ws.Range("I113").Value = 1179.4
$ws.Range("K113").Value = 1179.4
$ws.Range("M113").Value = 990.5999999999999

# Row 122
$ws.Range("H122").Value = 4581.1665
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 4581.1665
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 13743.4995
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -18643.4995

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 7073.724
$ws.Range("I7").Value = 4081.875
$ws.Range("K7").Value = 4081.875
$ws.Range("M7").Value = -3969.875

# Row 22
$ws.Range("H22").Value = 721.38464
$ws.Range("I22").Value = 660
$ws.Range("J22").Value = 859.5
$ws.Range("K22").Value = 660
$ws.Range("L22").Value = 859.5
$ws.Range("M22").Value = -365
$ws.Range("N22").Value = -1449.5

# Row 27
$ws.Range("H27").Value = 721.38464
$ws.Range("I27").Value = 660
$ws.Range("J27").Value = 859.5
$ws.Range("K27").Value = 660
$ws.Range("L27").Value = 859.5
$ws.Range("M27").Value = -553
$ws.Range("N27").Value = -1073.5

# Row 40
$ws.Range("H40").Value = 4697.25
$ws.Range("I40").Value = 4389
$ws.Range("J40").Value = 4741.2856
$ws.Range("K40").Value = 4389
$ws.Range("L40").Value = 4741.2856
$ws.Range("M40").Value = -4253
$ws.Range("N40").Value = -5013.2856

# Row 46
$ws.Range("H46").Value = 11895.659
$ws.Range("J46").Value = 1169.5952
$ws.Range("L46").Value = 1169.5952
$ws.Range("N46").Value = -1545.5952

# Row 126
$ws.Range("H126").Value = 7073.724
$ws.Range("I126").Value = 4081.875
$ws.Range("K126").Value = 12245.625
$ws.Range("M126").Value = -9775.625

# Row 139
$ws.Range("H139").Value = 89430
$ws.Range("I139").Value = 89430
$ws.Range("K139").Value = 89430
$ws.Range("M139").Value = -84290

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 469.5
$ws.Range("I122").Value = 344.5
$ws.Range("J122").Value = 532
$ws.Range("K122").Value = 1033.5
$ws.Range("L122").Value = 1596
$ws.Range("M122").Value = 1416.5
$ws.Range("N122").Value = -6496

# Row 132
$ws.Range("H132").Value = 1798.5
$ws.Range("I132").Value = 1798.5
$ws.Range("K132").Value = 5395.5
$ws.Range("M132").Value = -2865.5
